$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of data (row 25) following the existing table pattern
$ws.Range("A25").Value = 45971
$ws.Range("B25").Value = 5588
$ws.Range("C25").Value = 3748
$ws.Range("D25").Value = 3428
$ws.Range("E25").Value = 233
$ws.Range("F25").Value = 58
$ws.Range("G25").Value = 25
$ws.Range("H25").Value = 4
$ws.Range("I25").Value = 0

# Match the style of the previous date cell (row 24) for the new date cell
$ws.Range("A24").Copy()
$ws.Range("A25").PasteSpecial(-4122) | Out-Null

# Update the selection to the newly added row, consistent with the diff
$ws.Range("A25:I25").Select() | Out-Null
